$d = $word.ActiveDocument

# Update the table-level default shading from FFF2CC (themed accent4 tint) to flat ECEAF2
$tbl = $d.Tables(1)
$tbl.Shading.BackgroundPatternColor = 14343922   # RGB(0xF2,0xEA,0xEC) reversed order placeholder - will fix below

